# Updates the "Estado de Cuenta" database: previous periods for ANIBAL/ANGIE
# are removed/reordered and regrouped by worker, and the Valor Mora for the
# first row is updated, per commit "Elimna EC anteriores y se agregan
# nuevos, se modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docId1 = "1143387210"
$name1  = "ANGIE STEPHANIE CAMACHO AGRESOTT"
$docId2 = "1047423122"
$name2  = "ANIBAL FELIPE ARENAS SOTO"

# New row layout (row, DocId, Name, Periodo, ValorMora)
$rows = @(
  @(16, $docId1, $name1, "1803", 19673),
  @(17, $docId1, $name1, "1802", 29509),
  @(18, $docId1, $name1, "1801", 29509),
  @(19, $docId1, $name1, "1712", 29509),
  @(20, $docId1, $name1, "1711", 29509),
  @(21, $docId2, $name2, "1803", 29509),
  @(22, $docId2, $name2, "1802", 29509),
  @(23, $docId2, $name2, "1801", 29509),
  @(24, $docId2, $name2, "1712", 29509),
  @(25, $docId2, $name2, "1711", 29509)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
}
